$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number-format on Price cells whose new values parse as numbers,
# so Excel stores them as text (matching the source data's inlineStr type)
# instead of silently converting them to numeric cells.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply the updated Price / Volume(1h) values
$ws.Range("D2").Value = "34.359.58"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "1.802.28"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "227.27"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("E6").Value = "  +3.66%  "
$ws.Range("D8").Value = "35.85"
$ws.Range("E8").Value = "  +8.85%  "
$ws.Range("D9").Value = "0.301"
$ws.Range("E9").Value = "  +2.16%  "
$ws.Range("D10").Value = "0.0693"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("E11").Value = "  +1.98%  "
$ws.Range("D12").Value = "2.063.10"
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("E13").Value = "  +2.35%  "
$ws.Range("D14").Value = "1.801.88"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("E16").Value = "  +4.91%  "
$ws.Range("D17").Value = "34.367.84"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "69.08"
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("D19").Value = "245.65"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").Value = "11.47"
$ws.Range("E21").Value = "  +1.88%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "4.19"
$ws.Range("E23").Value = "  +0.89%  "
$ws.Range("E24").Value = "  +3.38%  "
$ws.Range("D25").Value = "170.82"
$ws.Range("E25").Value = "  +1.17%  "
$ws.Range("D26").Value = "7.88"
$ws.Range("E26").Value = "  +7.63%  "
$ws.Range("E27").Value = "  +2.09%  "
$ws.Range("E28").Value = "  +2.52%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").Value = "0.0531"
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("E32").Value = "  +1.17%  "
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("D35").Value = "1.394.00"
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("D37").Value = "2.47"
$ws.Range("E37").Value = "  -4.13%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("E40").Value = "  +11.63%  "
$ws.Range("D41").Value = "0.961"
$ws.Range("E41").Value = "  +2.50%  "
$ws.Range("D42").Value = "2.82"
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("D43").Value = "82.19"
$ws.Range("E43").Value = "  -2.66%  "
$ws.Range("D44").Value = "2.41"
$ws.Range("E44").Value = "  +0.34%  "
$ws.Range("D45").Value = "13.52"
$ws.Range("E45").Value = "  -3.61%  "
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("D47").Value = "0.0502"
$ws.Range("E47").Value = "  -5.04%  "
$ws.Range("D48").Value = "1.963.47"
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("D49").Value = "104.85"
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("E51").Value = "  +0.92%  "

Write-Output "Applied crypto price/volume updates"
